$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 used to hold "Total time taken for the ride" with a [hh]:mm:ss time
# style on B1. A new "Date and Time" text row is inserted above it, so that
# time style now belongs on B2 instead.
$ws.Range("B1").Style = "Normal"
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

$ws.Cells.Item(1, 1).Value = "Date and Time"
$ws.Cells.Item(1, 2).Value = "2024-03-12 10:20:08.465000 to 2024-03-12 11:29:19.295000"
$ws.Cells.Item(2, 1).Value = "Total time taken for the ride"
$ws.Cells.Item(2, 2).Value = 0.04805334490740741
$ws.Cells.Item(3, 1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(3, 2).Value = 34.70677833333333
$ws.Cells.Item(4, 1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(4, 2).Value = 1763.1729981775
$ws.Cells.Item(5, 1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(5, 2).Value = 3.652
$ws.Cells.Item(6, 1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(6, 2).Value = 4.698
$ws.Cells.Item(7, 1).Value = "Starting SoC (%)"
$ws.Cells.Item(7, 2).Value = 99
$ws.Cells.Item(8, 1).Value = "Ending SoC (%)"
$ws.Cells.Item(8, 2).Value = 9
$ws.Cells.Item(9, 1).Value = "Total distance covered (km)"
$ws.Cells.Item(9, 2).Value = 41.84414046820817
$ws.Cells.Item(10, 1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10, 2).Value = 42.13667620959026
$ws.Cells.Item(11, 1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(11, 2).Value = 90
$ws.Cells.Item(12, 1).Value = "Mode"
$ws.Cells.Item(12, 2).Value = "Custom mode`n69.09%`nEco mode`n26.05%`nSports mode`n0.17%"
$ws.Cells.Item(13, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(13, 2).Value = 6126.446865
$ws.Cells.Item(14, 1).Value = "Average Power(kW)"
$ws.Cells.Item(14, 2).Value = -1535.790658949674
$ws.Cells.Item(15, 1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(15, 2).Value = 108.2855966211111
$ws.Cells.Item(16, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(16, 2).Value = 5.786160427063244
$ws.Cells.Item(17, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 3.379
$ws.Cells.Item(18, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(18, 2).Value = 2.93
$ws.Cells.Item(19, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19, 2).Value = 0.4489999999999998
$ws.Cells.Item(20, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20, 2).Value = 28
$ws.Cells.Item(21, 1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(21, 2).Value = 46
$ws.Cells.Item(22, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(22, 2).Value = 18
$ws.Cells.Item(23, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23, 2).Value = 72
$ws.Cells.Item(24, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24, 2).Value = 70
$ws.Cells.Item(25, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25, 2).Value = 65
$ws.Cells.Item(26, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26, 2).Value = 49
$ws.Cells.Item(27, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(28, 1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(29, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(29, 2).Value = 46
$ws.Cells.Item(30, 1).Value = "lowest cell temp(C)"
$ws.Cells.Item(30, 2).Value = 26
$ws.Cells.Item(31, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(31, 2).Value = 20
$ws.Cells.Item(32, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(32, 2).Value = 55
$ws.Cells.Item(33, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(33, 2).Value = 1.908872808333333
$ws.Cells.Item(34, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(34, 2).Value = 0.0000001277384839217681
$ws.Cells.Item(35, 1).Value = "Cycle Count of battery"
$ws.Cells.Item(35, 2).Value = 137
$ws.Cells.Item(36, 1).Value = "Idling time percentage"
$ws.Cells.Item(36, 2).Value = 6.664304335720975
$ws.Cells.Item(37, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(37, 2).Value = 10.19008327216584
$ws.Cells.Item(38, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(38, 2).Value = 5.044420258675238
$ws.Cells.Item(39, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(39, 2).Value = 9.16500037966033
$ws.Cells.Item(40, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(40, 2).Value = 23.28836467640285
$ws.Cells.Item(41, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(41, 2).Value = 18.38315320559872
$ws.Cells.Item(42, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(42, 2).Value = 17.40616062163052
$ws.Cells.Item(43, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(43, 2).Value = 9.127034346604571
$ws.Cells.Item(44, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(44, 2).Value = 0.3999088815206662
$ws.Cells.Item(45, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(45, 2).Value = 0
